$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = "57.988.05"
$ws.Range("E2").Value = "  +2.63%  "

$ws.Range("D3").Value = "2.360.57"
$ws.Range("E3").Value = "  +1.27%  "

$ws.Range("D4").NumberFormat = "@"
$ws.Range("D4").Value = "0.997"
$ws.Range("D4").ClearFormats()
$ws.Range("E4").Value = "  -0.31%  "

$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = "541.73"
$ws.Range("D5").ClearFormats()
$ws.Range("E5").Value = "  +5.59%  "

$ws.Range("D6").NumberFormat = "@"
$ws.Range("D6").Value = "134.59"
$ws.Range("D6").ClearFormats()
$ws.Range("E6").Value = "  +1.88%  "

$ws.Range("D7").NumberFormat = "@"
$ws.Range("D7").Value = "0.998"
$ws.Range("D7").ClearFormats()
$ws.Range("E7").Value = "  -0.11%  "

$ws.Range("D8").NumberFormat = "@"
$ws.Range("D8").Value = "0.537"
$ws.Range("D8").ClearFormats()
$ws.Range("E8").Value = "  +0.63%  "

$ws.Range("D9").Value = "2.359.45"
$ws.Range("E9").Value = "  +1.03%  "

$ws.Range("E10").Value = "  +1.99%  "

$ws.Range("E11").Value = "  +0.94%  "

$ws.Range("D12").NumberFormat = "@"
$ws.Range("D12").Value = "5.39"
$ws.Range("D12").ClearFormats()
$ws.Range("E12").Value = "  +1.85%  "

$ws.Range("D13").NumberFormat = "@"
$ws.Range("D13").Value = "0.353"
$ws.Range("D13").ClearFormats()
$ws.Range("E13").Value = "  +4.20%  "

$ws.Range("D14").Value = "2.752.91"
$ws.Range("E14").Value = "  +0.22%  "

$ws.Range("D15").NumberFormat = "@"
$ws.Range("D15").Value = "23.51"
$ws.Range("D15").ClearFormats()
$ws.Range("E15").Value = "  -0.22%  "

$ws.Range("D16").Value = "57.910.75"
$ws.Range("E16").Value = "  +2.51%  "

$ws.Range("E17").Value = "  +0.87%  "

$ws.Range("D18").Value = "2.335.87"
$ws.Range("E18").Value = "  +0.45%  "

$ws.Range("E19").Value = "  +1.17%  "

$ws.Range("D20").NumberFormat = "@"
$ws.Range("D20").Value = "335.13"
$ws.Range("D20").ClearFormats()
$ws.Range("E20").Value = "  +3.25%  "

$ws.Range("E21").Value = "  +1.69%  "

$ws.Range("E22").Value = "  +1.82%  "

$ws.Range("D23").NumberFormat = "@"
$ws.Range("D23").Value = "0.998"
$ws.Range("D23").ClearFormats()
$ws.Range("E23").Value = "  +0.00%  "

$ws.Range("D24").NumberFormat = "@"
$ws.Range("D24").Value = "62.08"
$ws.Range("D24").ClearFormats()
$ws.Range("E24").Value = "  +0.37%  "

$ws.Range("E25").Value = "  +3.92%  "

$ws.Range("D26").NumberFormat = "@"
$ws.Range("D26").Value = "8.47"
$ws.Range("D26").ClearFormats()
$ws.Range("E26").Value = "  -2.41%  "

$ws.Range("D27").NumberFormat = "@"
$ws.Range("D27").Value = "0.995"
$ws.Range("D27").ClearFormats()
$ws.Range("E27").Value = "  -0.37%  "

$ws.Range("E28").Value = "  +7.87%  "

$ws.Range("E29").Value = "  +4.40%  "

$ws.Range("D30").NumberFormat = "@"
$ws.Range("D30").Value = "170.61"
$ws.Range("D30").ClearFormats()
$ws.Range("E30").Value = "  +1.74%  "

$ws.Range("E31").Value = "  +2.12%  "

$ws.Range("D32").NumberFormat = "@"
$ws.Range("D32").Value = "6.17"
$ws.Range("D32").ClearFormats()
$ws.Range("E32").Value = "  +0.95%  "

$ws.Range("E33").Value = "  +1.26%  "

$ws.Range("E34").Value = "  +15.24%  "

$ws.Range("D35").NumberFormat = "@"
$ws.Range("D35").Value = "0.999"
$ws.Range("D35").ClearFormats()
$ws.Range("E35").Value = "  -0.01%  "

$ws.Range("D36").NumberFormat = "@"
$ws.Range("D36").Value = "0.999"
$ws.Range("D36").ClearFormats()
$ws.Range("E36").Value = "  +0.19%  "

$ws.Range("E37").Value = "  -0.41%  "

$ws.Range("E38").Value = "  +4.88%  "

$ws.Range("E39").Value = "  +3.32%  "

$ws.Range("D40").NumberFormat = "@"
$ws.Range("D40").Value = "39.40"
$ws.Range("D40").ClearFormats()
$ws.Range("E40").Value = "  +2.47%  "

$ws.Range("D41").NumberFormat = "@"
$ws.Range("D41").Value = "149.76"
$ws.Range("D41").ClearFormats()
$ws.Range("E41").Value = "  -2.37%  "

$ws.Range("E42").Value = "  +0.59%  "

$ws.Range("E43").Value = "  +1.33%  "

$ws.Range("D44").NumberFormat = "@"
$ws.Range("D44").Value = "283.68"
$ws.Range("D44").ClearFormats()
$ws.Range("E44").Value = "  +1.61%  "

$ws.Range("D45").NumberFormat = "@"
$ws.Range("D45").Value = "19.28"
$ws.Range("D45").ClearFormats()
$ws.Range("E45").Value = "  +6.27%  "

$ws.Range("D46").NumberFormat = "@"
$ws.Range("D46").Value = "0.0931"
$ws.Range("D46").ClearFormats()
$ws.Range("E46").Value = "  +0.40%  "

$ws.Range("E47").Value = "  +1.94%  "

$ws.Range("E48").Value = "  +0.65%  "

$ws.Range("D49").NumberFormat = "@"
$ws.Range("D49").Value = "0.0218"
$ws.Range("D49").ClearFormats()
$ws.Range("E49").Value = "  +1.91%  "

$ws.Range("D50").NumberFormat = "@"
$ws.Range("D50").Value = "17.57"
$ws.Range("D50").ClearFormats()
$ws.Range("E50").Value = "  +2.51%  "

$ws.Range("E51").Value = "  +0.02%  "
